$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.1172255
$ws.Range("H2").Value = 4.234451
$ws.Range("I2").Value = 0.006295392006363395
$ws.Range("J2").Value = 0.004213039461358209
$ws.Range("M2").Value = 7.369448
$ws.Range("N2").Value = 14.738896
$ws.Range("O2").Value = 0.7452608427984224
$ws.Range("P2").Value = 0.661061693471796
$ws.Range("Q2").Value = 15.602783226524
$ws.Range("R2").Value = 62.411132906096
$ws.Range("S2").Value = 0.004691709152408835
$ws.Range("T2").Value = 0.00278507900098896
$ws.Range("G3").Value = 2.1172255
$ws.Range("H3").Value = 4.234451
$ws.Range("I3").Value = 0.006295392006363395
$ws.Range("J3").Value = 0.004213039461358209
$ws.Range("O3").Value = 0.01116592909756377
$ws.Range("P3").Value = 0.01485661309677453
$ws.Range("Q3").Value = 0.2337699248733333
$ws.Range("R3").Value = 1.40261954924
$ws.Range("S3").Value = 0.00007029390078442341
$ws.Range("T3").Value = 0.00006259149723884226
$ws.Range("G4").Value = 2.1172255
$ws.Range("H4").Value = 4.234451
$ws.Range("I4").Value = 0.006295392006363395
$ws.Range("J4").Value = 0.004213039461358209
$ws.Range("M4").Value = 0.084843
$ws.Range("N4").Value = 0.254529
$ws.Range("O4").Value = 0.008580040959044227
$ws.Range("P4").Value = 0.0114160091622658
$ws.Range("Q4").Value = 0.1796317630965
$ws.Range("R4").Value = 1.077790578579
$ws.Range("S4").Value = 0.00005401472126783754
$ws.Range("T4").Value = 0.00004809609709185269
$ws.Range("G5").Value = 2.1172255
$ws.Range("H5").Value = 4.234451
$ws.Range("I5").Value = 0.006295392006363395
$ws.Range("J5").Value = 0.004213039461358209
$ws.Range("M5").Value = 2.32371
$ws.Range("N5").Value = 6.97113
$ws.Range("O5").Value = 0.2349931871449696
$ws.Range("P5").Value = 0.3126656842691638
$ws.Range("Q5").Value = 4.919818066605
$ws.Range("R5").Value = 29.51890839963
$ws.Range("S5").Value = 0.001479374231902299
$ws.Range("T5").Value = 0.001317272866038554
$ws.Range("I6").Value = 0.004671287948366863
$ws.Range("J6").Value = 0.004689220411201648
$ws.Range("M6").Value = 7.369448
$ws.Range("N6").Value = 14.738896
$ws.Range("O6").Value = 0.7452608427984224
$ws.Range("P6").Value = 0.661061693471796
$ws.Range("Q6").Value = 11.57753054509867
$ws.Range("R6").Value = 69.46518327059201
$ws.Range("S6").Value = 0.003481327993354002
$ws.Range("T6").Value = 0.003099863986091473
$ws.Range("I7").Value = 0.004671287948366863
$ws.Range("J7").Value = 0.004689220411201648
$ws.Range("O7").Value = 0.01116592909756377
$ws.Range("P7").Value = 0.01485661309677453
$ws.Range("S7").Value = 0.00005215927002576853
$ws.Range("T7").Value = 0.00006966593337472085
$ws.Range("I8").Value = 0.004671287948366863
$ws.Range("J8").Value = 0.004689220411201648
$ws.Range("M8").Value = 0.084843
$ws.Range("N8").Value = 0.254529
$ws.Range("O8").Value = 0.008580040959044227
$ws.Range("P8").Value = 0.0114160091622658
$ws.Range("Q8").Value = 0.133289823612
$ws.Range("R8").Value = 1.199608412508
$ws.Range("S8").Value = 0.00004007984192847735
$ws.Range("T8").Value = 0.00005353218317816182
$ws.Range("I9").Value = 0.004671287948366863
$ws.Range("J9").Value = 0.004689220411201648
$ws.Range("M9").Value = 2.32371
$ws.Range("N9").Value = 6.97113
$ws.Range("O9").Value = 0.2349931871449696
$ws.Range("P9").Value = 0.3126656842691638
$ws.Range("Q9").Value = 3.65058868764
$ws.Range("R9").Value = 32.85529818876
$ws.Range("S9").Value = 0.001097720843058615
$ws.Range("T9").Value = 0.001466158308557293
$ws.Range("G10").Value = 75.40439600000001
$ws.Range("H10").Value = 226.213188
$ws.Range("I10").Value = 0.2242086314485916
$ws.Range("J10").Value = 0.2250693390296979
$ws.Range("M10").Value = 7.369448
$ws.Range("N10").Value = 14.738896
$ws.Range("O10").Value = 0.7452608427984224
$ws.Range("P10").Value = 0.661061693471796
$ws.Range("Q10").Value = 555.6887752934081
$ws.Range("R10").Value = 3334.132651760448
$ws.Range("S10").Value = 0.1670939136360582
$ws.Range("T10").Value = 0.1487847184075499
$ws.Range("G11").Value = 75.40439600000001
$ws.Range("H11").Value = 226.213188
$ws.Range("I11").Value = 0.2242086314485916
$ws.Range("J11").Value = 0.2250693390296979
$ws.Range("O11").Value = 0.01116592909756377
$ws.Range("P11").Value = 0.01485661309677453
$ws.Range("Q11").Value = 8.325650710346666
$ws.Range("R11").Value = 74.93085639312
$ws.Range("S11").Value = 0.00250349768181678
$ws.Range("T11").Value = 0.003343768089910996
$ws.Range("G12").Value = 75.40439600000001
$ws.Range("H12").Value = 226.213188
$ws.Range("I12").Value = 0.2242086314485916
$ws.Range("J12").Value = 0.2250693390296979
$ws.Range("M12").Value = 0.084843
$ws.Range("N12").Value = 0.254529
$ws.Range("O12").Value = 0.008580040959044227
$ws.Range("P12").Value = 0.0114160091622658
$ws.Range("Q12").Value = 6.397535169828001
$ws.Range("R12").Value = 57.577816528452
$ws.Range("S12").Value = 0.001923719241200167
$ws.Range("T12").Value = 0.002569393636508139
$ws.Range("G13").Value = 75.40439600000001
$ws.Range("H13").Value = 226.213188
$ws.Range("I13").Value = 0.2242086314485916
$ws.Range("J13").Value = 0.2250693390296979
$ws.Range("M13").Value = 2.32371
$ws.Range("N13").Value = 6.97113
$ws.Range("O13").Value = 0.2349931871449696
$ws.Range("P13").Value = 0.3126656842691638
$ws.Range("Q13").Value = 175.21794902916
$ws.Range("R13").Value = 1576.96154126244
$ws.Range("S13").Value = 0.0526875008895164
$ws.Range("T13").Value = 0.07037145889572891
$ws.Range("G14").Value = 1.7411535
$ws.Range("H14").Value = 3.482307
$ws.Range("I14").Value = 0.005177173534775417
$ws.Range("J14").Value = 0.003464698684094803
$ws.Range("M14").Value = 7.369448
$ws.Range("N14").Value = 14.738896
$ws.Range("O14").Value = 0.7452608427984224
$ws.Range("P14").Value = 0.661061693471796
$ws.Range("Q14").Value = 12.831340178268
$ws.Range("R14").Value = 51.32536071307201
$ws.Range("S14").Value = 0.003858344711840415
$ws.Range("T14").Value = 0.002290379579477214
$ws.Range("G15").Value = 1.7411535
$ws.Range("H15").Value = 3.482307
$ws.Range("I15").Value = 0.005177173534775417
$ws.Range("J15").Value = 0.003464698684094803
$ws.Range("O15").Value = 0.01116592909756377
$ws.Range("P15").Value = 0.01485661309677453
$ws.Range("Q15").Value = 0.19224656178
$ws.Range("R15").Value = 1.15347937068
$ws.Range("S15").Value = 0.00005780795261508591
$ws.Range("T15").Value = 0.00005147368784650032
$ws.Range("G16").Value = 1.7411535
$ws.Range("H16").Value = 3.482307
$ws.Range("I16").Value = 0.005177173534775417
$ws.Range("J16").Value = 0.003464698684094803
$ws.Range("M16").Value = 0.084843
$ws.Range("N16").Value = 0.254529
$ws.Range("O16").Value = 0.008580040959044227
$ws.Range("P16").Value = 0.0114160091622658
$ws.Range("Q16").Value = 0.1477246864005
$ws.Range("R16").Value = 0.8863481184030001
$ws.Range("S16").Value = 0.00004442036098045286
$ws.Range("T16").Value = 0.00003955303192211653
$ws.Range("G17").Value = 1.7411535
$ws.Range("H17").Value = 3.482307
$ws.Range("I17").Value = 0.005177173534775417
$ws.Range("J17").Value = 0.003464698684094803
$ws.Range("M17").Value = 2.32371
$ws.Range("N17").Value = 6.97113
$ws.Range("O17").Value = 0.2349931871449696
$ws.Range("P17").Value = 0.3126656842691638
$ws.Range("Q17").Value = 4.045935799485001
$ws.Range("R17").Value = 24.27561479691001
$ws.Range("S17").Value = 0.001216600509339464
$ws.Range("T17").Value = 0.001083292384848973
$ws.Range("G18").Value = 197.2895866666667
$ws.Range("H18").Value = 591.86876
$ws.Range("I18").Value = 0.586623997698909
$ws.Range("J18").Value = 0.5888759704209946
$ws.Range("M18").Value = 7.369448
$ws.Range("N18").Value = 14.738896
$ws.Range("O18").Value = 0.7452608427984224
$ws.Range("P18").Value = 0.661061693471796
$ws.Range("Q18").Value = 1453.915349881493
$ws.Range("R18").Value = 8723.49209928896
$ws.Range("S18").Value = 0.4371878949308688
$ws.Range("T18").Value = 0.3892833462513499
$ws.Range("G19").Value = 197.2895866666667
$ws.Range("H19").Value = 591.86876
$ws.Range("I19").Value = 0.586623997698909
$ws.Range("J19").Value = 0.5888759704209946
$ws.Range("O19").Value = 0.01116592909756377
$ws.Range("P19").Value = 0.01485661309677453
$ws.Range("Q19").Value = 21.78340089582222
$ws.Range("R19").Value = 196.0506080624
$ws.Range("S19").Value = 0.006550201965235431
$ws.Range("T19").Value = 0.008748702454532357
$ws.Range("G20").Value = 197.2895866666667
$ws.Range("H20").Value = 591.86876
$ws.Range("I20").Value = 0.586623997698909
$ws.Range("J20").Value = 0.5888759704209946
$ws.Range("M20").Value = 0.084843
$ws.Range("N20").Value = 0.254529
$ws.Range("O20").Value = 0.008580040959044227
$ws.Range("P20").Value = 0.0114160091622658
$ws.Range("Q20").Value = 16.73864040156
$ws.Range("R20").Value = 150.64776361404
$ws.Range("S20").Value = 0.005033257927814906
$ws.Range("T20").Value = 0.006722613473764239
$ws.Range("G21").Value = 197.2895866666667
$ws.Range("H21").Value = 591.86876
$ws.Range("I21").Value = 0.586623997698909
$ws.Range("J21").Value = 0.5888759704209946
$ws.Range("M21").Value = 2.32371
$ws.Range("N21").Value = 6.97113
$ws.Range("O21").Value = 0.2349931871449696
$ws.Range("P21").Value = 0.3126656842691638
$ws.Range("Q21").Value = 458.4437854332
$ws.Range("R21").Value = 4125.9940688988
$ws.Range("S21").Value = 0.13785264287499
$ws.Range("T21").Value = 0.1841213082413481
$ws.Range("G22").Value = 58.19014966666668
$ws.Range("H22").Value = 174.570449
$ws.Range("I22").Value = 0.1730235173629937
$ws.Range("J22").Value = 0.1736877319926528
$ws.Range("M22").Value = 7.369448
$ws.Range("N22").Value = 14.738896
$ws.Range("O22").Value = 0.7452608427984224
$ws.Range("P22").Value = 0.661061693471796
$ws.Range("Q22").Value = 428.8292820807174
$ws.Range("R22").Value = 2572.975692484305
$ws.Range("S22").Value = 0.1289476523738921
$ws.Range("T22").Value = 0.1148183062463385
$ws.Range("G23").Value = 58.19014966666668
$ws.Range("H23").Value = 174.570449
$ws.Range("I23").Value = 0.1730235173629937
$ws.Range("J23").Value = 0.1736877319926528
$ws.Range("O23").Value = 0.01116592909756377
$ws.Range("P23").Value = 0.01485661309677453
$ws.Range("Q23").Value = 6.424968391862222
$ws.Range("R23").Value = 57.82471552676
$ws.Range("S23").Value = 0.001931968327086282
$ws.Range("T23").Value = 0.00258041143387111
$ws.Range("G24").Value = 58.19014966666668
$ws.Range("H24").Value = 174.570449
$ws.Range("I24").Value = 0.1730235173629937
$ws.Range("J24").Value = 0.1736877319926528
$ws.Range("M24").Value = 0.084843
$ws.Range("N24").Value = 0.254529
$ws.Range("O24").Value = 0.008580040959044227
$ws.Range("P24").Value = 0.0114160091622658
$ws.Range("Q24").Value = 4.937026868169001
$ws.Range("R24").Value = 44.43324181352101
$ws.Range("S24").Value = 0.001484548865852386
$ws.Range("T24").Value = 0.001982820739801292
$ws.Range("G25").Value = 58.19014966666668
$ws.Range("H25").Value = 174.570449
$ws.Range("I25").Value = 0.1730235173629937
$ws.Range("J25").Value = 0.1736877319926528
$ws.Range("M25").Value = 2.32371
$ws.Range("N25").Value = 6.97113
$ws.Range("O25").Value = 0.2349931871449696
$ws.Range("P25").Value = 0.3126656842691638
$ws.Range("Q25").Value = 135.21703268193
$ws.Range("R25").Value = 1216.95329413737
$ws.Range("S25").Value = 0.04065934779616288
$ws.Range("T25").Value = 0.05430619357264193

Write-Output "Applied 266 cell updates"
